$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style-only changes (copy the existing "no border-fill-but-smaller-font" style from C60
# onto C19 and C71 before C60 itself is edited) ---
$ws.Cells.Item(60,3).Copy()
$ws.Cells.Item(19,3).PasteSpecial(-4122)
$ws.Cells.Item(71,3).PasteSpecial(-4122)

# --- C60 needs the normal "s=12" cell style (copy it from a cell that keeps that style, e.g. C4) ---
$ws.Cells.Item(4,3).Copy()
$ws.Cells.Item(60,3).PasteSpecial(-4122)

# --- Numeric FP values ---
$ws.Cells.Item(2,3).Value = 31914
$ws.Cells.Item(5,3).Value = 92322
$ws.Cells.Item(8,3).Value = 31914
$ws.Cells.Item(10,3).Value = 31914
$ws.Cells.Item(13,3).Value = 92322
$ws.Cells.Item(20,3).Value = 31914
$ws.Cells.Item(26,3).Value = 92322
$ws.Cells.Item(45,3).Value = 9999
$ws.Cells.Item(46,3).Value = 92322
$ws.Cells.Item(47,3).Value = 92322
$ws.Cells.Item(50,3).Value = 9999
$ws.Cells.Item(51,3).Value = 9999
$ws.Cells.Item(52,3).Value = 14984
$ws.Cells.Item(58,3).Value = 31914
$ws.Cells.Item(59,3).Value = 92322
$ws.Cells.Item(60,3).Value = 31710
$ws.Cells.Item(61,3).Value = 9999
$ws.Cells.Item(62,3).Value = 9999
$ws.Cells.Item(63,3).Value = 9999
$ws.Cells.Item(64,3).Value = 9999
$ws.Cells.Item(68,3).Value = 92322
$ws.Cells.Item(70,3).Value = 9999
$ws.Cells.Item(79,3).Value = 56497
$ws.Cells.Item(84,3).Value = 31914
$ws.Cells.Item(89,3).Value = 14984
$ws.Cells.Item(90,3).Value = 14984
$ws.Cells.Item(110,3).Value = 9999
$ws.Cells.Item(111,3).Value = 92322
$ws.Cells.Item(123,3).Value = 56610

# --- String FP values (matched/created as shared strings automatically) ---
$ws.Cells.Item(3,3).Value = "14978, 9999"
$ws.Cells.Item(12,3).Value = "92322, 33836"
$ws.Cells.Item(36,3).Value = "31914, 55873"
$ws.Cells.Item(41,3).Value = "9999, 92322, 33836"
$ws.Cells.Item(42,3).Value = "92322, 9999"
$ws.Cells.Item(48,3).Value = "56610, 9999"
$ws.Cells.Item(49,3).Value = "56610, 9999"
$ws.Cells.Item(56,3).Value = "53310, 94328"
$ws.Cells.Item(66,3).Value = "93031, 9999, 14984, 14978"
$ws.Cells.Item(67,3).Value = "30520, 9999"
$ws.Cells.Item(76,3).Value = "14984, 9999"
$ws.Cells.Item(95,3).Value = "93031, 9999, 14978"
$ws.Cells.Item(96,3).Value = "9999, 31914"
$ws.Cells.Item(109,3).Value = "93031, 9999"
$ws.Cells.Item(119,3).Value = "93031, 9999"
$ws.Cells.Item(122,3).Value = "35364, 9999, 30520, 59582"
$ws.Cells.Item(126,3).Value = "9999, 56610"
$ws.Cells.Item(132,3).Value = "14978, 9999, 14984"
$ws.Cells.Item(134,3).Value = "14984, 9999"
$ws.Cells.Item(136,3).Value = "9999, 14984"

# --- Row heights that changed because of the new wrapped FP text ---
$ws.Rows.Item(12).RowHeight = 53
$ws.Rows.Item(36).RowHeight = 53
$ws.Rows.Item(41).RowHeight = 70.5
$ws.Rows.Item(56).RowHeight = 53
$ws.Rows.Item(66).RowHeight = 105.5
$ws.Rows.Item(67).RowHeight = 35.5

# --- New column width for column C ---
$ws.Columns.Item(3).ColumnWidth = 30.6

# --- Selection moves to F5 ---
$ws.Range("F5").Select()
